$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Worker roster (Tipo Doc stays "CC" for all, unchanged)
$workers = @(
    @("73159284", "WILMAN MALLARINO MARTINEZ"),
    @("73191182", "JAINER GOMEZ MARIMON"),
    @("1047368039", "ELKIN ZUÑIGA QUINTANA"),
    @("73098243", "WILSON MARTINEZ VERGARA"),
    @("15049482", "JUAN MIGUEL CALLE BEDOYA"),
    @("73184551", "CIRIN ESTEBAN CABARCAS PUERTAS")
)

$periods = @("1706", "1707", "1708")

# Data table previously was grouped by worker (3 periods each); now it is
# regrouped by period (6 workers each), covering rows 16-33.
$row = 16
foreach ($period in $periods) {
    foreach ($worker in $workers) {
        $ws.Cells.Item($row, 3).Value = $worker[0]
        $ws.Cells.Item($row, 4).Value = $worker[1]
        $ws.Cells.Item($row, 5).Value = $period
        $row = $row + 1
    }
}
